# Refactor code to save results in a specified folder
# Updates sliding window results (window size 10) with recomputed
# IPC PO / DELTA / DELTA^2 values after saving results to the new folder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 30.06443033626282
$ws.Cells.Item(2, 4).Value = 0.1444303362628148
$ws.Cells.Item(2, 5).Value = 0.02086012203298974
$ws.Cells.Item(3, 3).Value = 29.82430367137924
$ws.Cells.Item(3, 4).Value = -0.1556963286207598
$ws.Cells.Item(3, 5).Value = 0.02424134674598362
$ws.Cells.Item(4, 3).Value = 29.74799945601933
$ws.Cells.Item(4, 4).Value = -0.2920005439806665
$ws.Cells.Item(4, 5).Value = 0.08526431768500518
$ws.Cells.Item(5, 3).Value = 29.11472861437268
$ws.Cells.Item(5, 4).Value = -1.095271385627317
$ws.Cells.Item(5, 5).Value = 1.199619408173983
$ws.Cells.Item(6, 3).Value = 29.52318149212922
$ws.Cells.Item(6, 4).Value = -0.6968185078707805
$ws.Cells.Item(6, 5).Value = 0.485556032911261
$ws.Cells.Item(7, 3).Value = 29.77955219312792
$ws.Cells.Item(7, 4).Value = -0.6004478068720829
$ws.Cells.Item(7, 5).Value = 0.3605375687774942
$ws.Cells.Item(8, 3).Value = 30.23222788438841
$ws.Cells.Item(8, 4).Value = -0.2077721156115935
$ws.Cells.Item(8, 5).Value = 0.04316925202571736
$ws.Cells.Item(9, 3).Value = 30.49980616390415
$ws.Cells.Item(9, 4).Value = 0.01980616390414625
$ws.Cells.Item(9, 5).Value = 0.0003922841285979057
$ws.Cells.Item(10, 3).Value = 30.54786367970804
$ws.Cells.Item(10, 4).Value = -0.1421363202919572
$ws.Cells.Item(10, 5).Value = 0.02020273354613785
$ws.Cells.Item(11, 3).Value = 30.96101268311433
$ws.Cells.Item(11, 4).Value = 0.2110126831143262
$ws.Cells.Item(11, 5).Value = 0.04452635243510705
$ws.Cells.Item(12, 3).Value = 31.31503111168806
$ws.Cells.Item(12, 4).Value = 0.375031111688056
$ws.Cells.Item(12, 5).Value = 0.1406483347339792
$ws.Cells.Item(13, 3).Value = 31.49457366232175
$ws.Cells.Item(13, 4).Value = 0.5445736623217527
$ws.Cells.Item(13, 5).Value = 0.2965604736945264
$ws.Cells.Item(14, 3).Value = 31.53957830086107
$ws.Cells.Item(14, 4).Value = 0.5195783008610704
$ws.Cells.Item(14, 5).Value = 0.269961610725677
$ws.Cells.Item(15, 3).Value = 31.54807976351414
$ws.Cells.Item(15, 4).Value = 0.4280797635141411
$ws.Cells.Item(15, 5).Value = 0.183252283930323
$ws.Cells.Item(16, 3).Value = 31.82866309169452
$ws.Cells.Item(16, 4).Value = 0.5486630916945145
$ws.Cells.Item(16, 5).Value = 0.3010311881877833
$ws.Cells.Item(17, 3).Value = 31.6646103044787
$ws.Cells.Item(17, 4).Value = 0.2846103044787007
$ws.Cells.Item(17, 5).Value = 0.08100302541545873
$ws.Cells.Item(18, 3).Value = 31.59493988768019
$ws.Cells.Item(18, 4).Value = 0.01493988768018895
$ws.Cells.Item(18, 5).Value = 0.0002232002438966614
$ws.Cells.Item(19, 3).Value = 31.68670405269421
$ws.Cells.Item(19, 4).Value = 0.03670405269421551
$ws.Cells.Item(19, 5).Value = 0.001347187484179749
$ws.Cells.Item(20, 3).Value = 32.21530595261872
$ws.Cells.Item(20, 4).Value = 0.3353059526187216
$ws.Cells.Item(20, 5).Value = 0.1124300818615484
$ws.Cells.Item(21, 3).Value = 32.17015382033146
$ws.Cells.Item(21, 4).Value = -0.109846179668537
$ws.Cells.Item(21, 5).Value = 0.01206618318777252
$ws.Cells.Item(22, 3).Value = 32.2114913262569
$ws.Cells.Item(22, 4).Value = -0.2385086737431052
$ws.Cells.Item(22, 5).Value = 0.05688638745069502
$ws.Cells.Item(23, 3).Value = 33.10141974170725
$ws.Cells.Item(23, 4).Value = 0.2514197417072523
$ws.Cells.Item(23, 5).Value = 0.06321188652014148
$ws.Cells.Item(24, 3).Value = 33.18490414560242
$ws.Cells.Item(24, 4).Value = 0.2849041456024253
$ws.Cells.Item(24, 5).Value = 0.08117037218144794
$ws.Cells.Item(25, 3).Value = 33.06751199284937
$ws.Cells.Item(25, 4).Value = -0.03248800715063282
$ws.Cells.Item(25, 5).Value = 0.00105547060861957
$ws.Cells.Item(26, 3).Value = 33.40909883880168
$ws.Cells.Item(26, 4).Value = 0.009098838801683939
$ws.Cells.Item(26, 5).Value = 0.00008278886753902923
$ws.Cells.Item(27, 3).Value = 33.62261684802064
$ws.Cells.Item(27, 4).Value = -0.07738315197936174
$ws.Cells.Item(27, 5).Value = 0.005988152210260997
$ws.Cells.Item(28, 3).Value = 34.10825457068023
$ws.Cells.Item(28, 4).Value = 0.008254570680229278
$ws.Cells.Item(28, 5).Value = 0.00006813793711490086
$ws.Cells.Item(29, 3).Value = 34.5251974927726
$ws.Cells.Item(29, 4).Value = 0.1251974927725996
$ws.Cells.Item(29, 5).Value = 0.01567441219654513
$ws.Cells.Item(30, 3).Value = 34.78657990434206
$ws.Cells.Item(30, 4).Value = -0.1134200956579434
$ws.Cells.Item(30, 5).Value = 0.01286411809905704
$ws.Cells.Item(31, 3).Value = 35.30601485221275
$ws.Cells.Item(31, 4).Value = 0.006014852212750554
$ws.Cells.Item(31, 5).Value = 0.00003617844714123024
$ws.Cells.Item(32, 3).Value = 35.48539191447748
$ws.Cells.Item(32, 4).Value = -0.2146080855225208
$ws.Cells.Item(32, 5).Value = 0.04605663037164158
$ws.Cells.Item(33, 3).Value = 35.67686639017901
$ws.Cells.Item(33, 4).Value = -0.6231336098209894
$ws.Cells.Item(33, 5).Value = 0.3882954956885371
$ws.Cells.Item(34, 3).Value = 36.22003163832851
$ws.Cells.Item(34, 4).Value = -0.5799683616714901
$ws.Cells.Item(34, 5).Value = 0.3363633005399124
$ws.Cells.Item(35, 3).Value = 36.82715857126485
$ws.Cells.Item(35, 4).Value = -0.4728414287351512
$ws.Cells.Item(35, 5).Value = 0.2235790167282991
$ws.Cells.Item(36, 3).Value = 37.80912703771192
$ws.Cells.Item(36, 4).Value = -0.09087296228808128
$ws.Cells.Item(36, 5).Value = 0.008257895275011043
$ws.Cells.Item(37, 3).Value = 38.5326960693319
$ws.Cells.Item(37, 4).Value = 0.03269606933189806
$ws.Cells.Item(37, 5).Value = 0.001069032949756285
$ws.Cells.Item(38, 3).Value = 39.31931509925381
$ws.Cells.Item(38, 4).Value = 0.4193150992538079
$ws.Cells.Item(38, 5).Value = 0.1758251524622308
$ws.Cells.Item(39, 3).Value = 39.84451746815164
$ws.Cells.Item(39, 4).Value = 0.4445174681516448
$ws.Cells.Item(39, 5).Value = 0.1975957794919486
$ws.Cells.Item(40, 3).Value = 40.13033370104905
$ws.Cells.Item(40, 4).Value = 0.2303337010490552
$ws.Cells.Item(40, 5).Value = 0.05305361383895551
$ws.Cells.Item(41, 3).Value = 39.87095759650498
$ws.Cells.Item(41, 4).Value = -0.229042403495022
$ws.Cells.Item(41, 5).Value = 0.05246042259877648
$ws.Cells.Item(42, 3).Value = 40.60442235086847
$ws.Cells.Item(42, 4).Value = 0.004422350868466651
$ws.Cells.Item(42, 5).Value = 0.00001955718720382775
$ws.Cells.Item(43, 3).Value = 40.77083204976728
$ws.Cells.Item(43, 4).Value = -0.1291679502327199
$ws.Cells.Item(43, 5).Value = 0.0166843593673224
$ws.Cells.Item(44, 3).Value = 41.48470126378374
$ws.Cells.Item(44, 4).Value = 0.2847012637837381
$ws.Cells.Item(44, 5).Value = 0.08105480960005759
$ws.Cells.Item(45, 3).Value = 41.10928511946057
$ws.Cells.Item(45, 4).Value = -0.3907148805394272
$ws.Cells.Item(45, 5).Value = 0.1526581178749388
$ws.Cells.Item(46, 3).Value = 41.28583983217257
$ws.Cells.Item(46, 4).Value = -0.514160167827427
$ws.Cells.Item(46, 5).Value = 0.2643606781803279
$ws.Cells.Item(47, 3).Value = 42.20022596496717
$ws.Cells.Item(47, 4).Value = 0.0002259649671714214
$ws.Cells.Item(47, 5).Value = 0.00000005106016638878157
$ws.Cells.Item(48, 3).Value = 43.26718534252096
$ws.Cells.Item(48, 4).Value = 0.5671853425209576
$ws.Cells.Item(48, 5).Value = 0.321699212770616
$ws.Cells.Item(49, 3).Value = 43.89964921480373
$ws.Cells.Item(49, 4).Value = 0.1996492148037277
$ws.Cells.Item(49, 5).Value = 0.03985980897174499
$ws.Cells.Item(50, 3).Value = 43.64463790310828
$ws.Cells.Item(50, 4).Value = -0.5553620968917272
$ws.Cells.Item(50, 5).Value = 0.3084270586639762
$ws.Cells.Item(51, 3).Value = 45.55437905176503
$ws.Cells.Item(51, 4).Value = -0.04562094823497631
$ws.Cells.Item(51, 5).Value = 0.002081270917858388

$ws.Cells.Item(52, 3).Value = -1.276610584994213
$ws.Cells.Item(52, 5).Value = 6.589332156985265
$ws.Cells.Item(53, 5).Value = 0.1317866431397053

Write-Host "Updated sliding window results (window size 10)"
